$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / rich-text header updates ---
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Cells changing from numeric to text placeholders ---
# Donor cell C14 already has the General-format "text" style used for placeholders
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D14").PasteSpecial(-4122)

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# --- Pure numeric value updates ---
$ws.Range("N15").Value = -58.333333333333
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -28.571428571428
$ws.Range("J16").Value = 103
$ws.Range("K16").Value = 65.04854368932
$ws.Range("L16").Value = 44.067796610169
$ws.Range("M16").Value = -33.852140077821
$ws.Range("N16").Value = -83.284169124877
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -36.363636363636
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -3.703703703703
$ws.Range("I17").Value = 357
$ws.Range("J17").Value = 336
$ws.Range("K17").Value = 6.25
$ws.Range("L17").Value = 66.822429906542
$ws.Range("M17").Value = 74.146341463414
$ws.Range("N17").Value = -49.718309859154
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = -12.5
$ws.Range("I18").Value = 156
$ws.Range("J18").Value = 165
$ws.Range("K18").Value = -5.454545454545
$ws.Range("L18").Value = 13.043478260869
$ws.Range("M18").Value = 0.645161290322
$ws.Range("N18").Value = -83.386581469648
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -38.095238095238
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -19.607843137254
$ws.Range("I19").Value = 494
$ws.Range("J19").Value = 403
$ws.Range("K19").Value = 22.58064516129
$ws.Range("L19").Value = 58.842443729903
$ws.Range("M19").Value = -6.439393939393
$ws.Range("N19").Value = -24.924012158054
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 115
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = 66.666666666666
$ws.Range("L20").Value = 79.6875
$ws.Range("M20").Value = 6.481481481481
$ws.Range("N20").Value = -89.722966934763
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -34.210526315789
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = -10.344827586206
$ws.Range("I21").Value = 1312
$ws.Range("J21").Value = 1096
$ws.Range("K21").Value = 19.70802919708
$ws.Range("L21").Value = 49.771689497716
$ws.Range("M21").Value = 2.740798747063
$ws.Range("N21").Value = -70.837963991998
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 92.857142857142
$ws.Range("L22").Value = -12.903225806451
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 12.5
$ws.Range("I23").Value = 172
$ws.Range("J23").Value = 176
$ws.Range("K23").Value = -2.272727272727
$ws.Range("L23").Value = 57.798165137614
$ws.Range("M23").Value = 48.275862068965
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 64.705882352941
$ws.Range("F24").Value = 82
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 12.328767123287
$ws.Range("I24").Value = 1160
$ws.Range("J24").Value = 833
$ws.Range("K24").Value = 39.255702280912
$ws.Range("L24").Value = 19.096509240246
$ws.Range("M24").Value = 2.112676056338
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -27.272727272727
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -23.91304347826
$ws.Range("I25").Value = 531
$ws.Range("J25").Value = 441
$ws.Range("K25").Value = 20.408163265306
$ws.Range("L25").Value = 32.089552238806
$ws.Range("M25").Value = 6.626506024096
$ws.Range("L26").Value = -6.666666666666
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 68.75
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 6
$ws.Range("J28").Value = 18
$ws.Range("K28").Value = 44.444444444444
$ws.Range("D29").Value = 1
$ws.Range("J29").Value = 15
$ws.Range("K29").Value = 13.333333333333